# Atualização automática dos dados: Wed Jan 21 09:38:00 UTC 2026
#
# Applies the refreshed "Entrada" dashboard figures. Row 2 and row 3 swap
# their category labels (FERRAMENTAS/ MATRIZARIA moves up to row 2,
# DEVOLUÇÃO moves down to row 3) along with the rest of that row's figures;
# the remaining rows keep their labels and only refresh the numeric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# The mojibake "DEVOLUÇÃO" text (double-UTF8-encoded in the source file) is
# reconstructed from its exact code points to avoid any re-encoding drift,
# for use on row 3 below (it moves there from row 2).
$devolucaoCodes = 0x0044,0x0045,0x0056,0x004f,0x004c,0x0055,0x00c3,0x0087,0x00c3,0x0083,0x004f
$devolucaoSb = New-Object System.Text.StringBuilder
foreach ($code in $devolucaoCodes) {
    [void]$devolucaoSb.Append([char]$code)
}
$devolucao = $devolucaoSb.ToString()

# --- Row 2 (was DEVOLUÇÃO) now becomes FERRAMENTAS/ MATRIZARIA ---
$ws.Range("A2").Value = "FERRAMENTAS/ MATRIZARIA"
$ws.Range("B2").Value = "R$ 447.144,85"
$ws.Range("D2").Value = "R$ 447.144,85"
$ws.Range("E2").Value = "R$ 680.000,00"
$ws.Range("F2").Value = "65,76 %"

# --- Row 3 (was FERRAMENTAS/ MATRIZARIA) now becomes DEVOLUÇÃO ---
$ws.Range("A3").Value = $devolucao
$ws.Range("B3").Value = "R$ 438.602,92"
$ws.Range("D3").Value = "R$ 438.602,92"
$ws.Range("E3").Value = "R$ 438.602,92"
$ws.Range("F3").Value = "100,00 %"

# --- Row 4 (MATERIA PRIMA) ---
$ws.Range("B4").Value = "R$ 368.181,47"
$ws.Range("D4").Value = "R$ 368.181,47"
$ws.Range("F4").Value = "36.818.147,00 %"

# --- Row 5 (REFUGO REAL (PROCESSO)) ---
$ws.Range("B5").Value = "R$ 283.784,61"
$ws.Range("D5").Value = "R$ 283.784,61"
$ws.Range("E5").Value = "R$ 283.784,61"

# --- Row 6 (was FRETES) now becomes MANUTENCAO ---
$ws.Range("A6").Value = "MANUTENCAO"
$ws.Range("B6").Value = "R$ 205.025,51"
$ws.Range("C6").Value = "R$ 253.396,69"
$ws.Range("D6").Value = "R$ 458.422,20"
$ws.Range("E6").Value = "R$ 480.000,00"
$ws.Range("F6").Value = "95,50 %"

# --- Row 7 (was MANUTENCAO) now becomes FRETES ---
$ws.Range("A7").Value = "FRETES"
$ws.Range("B7").Value = "R$ 185.370,49"
$ws.Range("C7").Value = "R$ 0,00"
$ws.Range("D7").Value = "R$ 185.370,49"
$ws.Range("E7").Value = "R$ 376.000,00"
$ws.Range("F7").Value = "49,30 %"

# --- Row 8 (REFUGO MP+CP*) ---
$ws.Range("B8").Value = "R$ 166.448,51"
$ws.Range("D8").Value = "R$ 166.448,51"
$ws.Range("E8").Value = "R$ 280.000,00"
$ws.Range("F8").Value = "59,45 %"

# --- Row 9 (CUSTO DESENVOLVIMENTO) ---
$ws.Range("B9").Value = "R$ 154.317,23"
$ws.Range("D9").Value = "R$ 154.317,23"
$ws.Range("E9").Value = "R$ 154.317,23"

# --- Row 10 (OLEOS E LUBRIFICANTES) ---
$ws.Range("C10").Value = "R$ 108.767,66"
$ws.Range("D10").Value = "R$ 190.859,64"
$ws.Range("F10").Value = "68,16 %"

# --- Row 11 (DESP. INDUSTRIAL) ---
$ws.Range("B11").Value = "R$ 59.268,75"
$ws.Range("C11").Value = "R$ 104.851,11"
$ws.Range("D11").Value = "R$ 164.119,86"
$ws.Range("E11").Value = "R$ 450.000,00"
$ws.Range("F11").Value = "36,47 %"

# --- Row 12 (EMBALAGENS) ---
$ws.Range("C12").Value = "R$ 72.265,01"
$ws.Range("D12").Value = "R$ 120.239,09"
$ws.Range("F12").Value = "75,15 %"

# --- Row 16 (MATERIAL QUALIDADE) ---
$ws.Range("C16").Value = "R$ 4.532,77"
$ws.Range("D16").Value = "R$ 4.643,14"
$ws.Range("F16").Value = "10,32 %"

# --- Row 18 (Total Geral) ---
$ws.Range("B18").Value = "R$ 2.332.555,21"
$ws.Range("C18").Value = "R$ 543.813,25"
$ws.Range("D18").Value = "R$ 2.876.368,46"
$ws.Range("E18").Value = "R$ 3.955.006,84"
$ws.Range("F18").Value = "72,73 %"
